# Batch19TestData1.xlsx edit script
# Rewrites the sample "registration form" rows with a new roster of names
# and appends a brand-new 8th row, matching the locally-committed changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 2: Jane/XHIXHI/Lopez -> Lil/Wayne/Rapper -------------------------
$ws.Range("A2").Value = "Lil"
$ws.Range("B2").Value = "Wayne"
$ws.Range("C2").Value = "Rapper"
$ws.Range("D2").Value = "/Users/shpendpllana/Desktop/PANDA.jpeg"
$ws.Range("E2").Value = "lilTheBest"

# ---- Row 3: John/JS/Smith -> John/Jones/Goat -------------------------------
$ws.Range("B3").Value = "Jones"
$ws.Range("C3").Value = "Goat"
$ws.Range("D3").Value = "/Users/shpendpllana/Desktop/PANDA.jpeg"
$ws.Range("E3").Value = "bones123"

# ---- Row 4: Laura/LR/Ricarso -> Eminem/lyricalist/emcylopedi ---------------
$ws.Range("A4").Value = "Eminem"
$ws.Range("B4").Value = "lyricalist"
$ws.Range("C4").Value = "emcylopedi"
$ws.Range("D4").Value = "/Users/shpendpllana/Desktop/PANDA.jpeg"
$ws.Range("E4").Value = "m&m"

# ---- Row 5: Sarah/SRB/Brown -> 50cent/gangsta/realG ------------------------
$ws.Range("A5").Value = "50cent"
$ws.Range("B5").Value = "gangsta"
$ws.Range("C5").Value = "realG"
$ws.Range("D5").Value = "/Users/shpendpllana/Desktop/PANDA.jpeg"
$ws.Range("E5").Value = "50gold"

# ---- Row 6: Linda/LCH/Christos -> SDET/QA/Career ---------------------------
$ws.Range("A6").Value = "SDET"
$ws.Range("B6").Value = "QA"
$ws.Range("C6").Value = "Career"
$ws.Range("D6").Value = "/Users/shpendpllana/Desktop/PANDA.jpeg"
$ws.Range("E6").Value = "debugger95"

# ---- Row 7: Joe/JP/Partiz -> Software/Tester/SDET --------------------------
$ws.Range("A7").Value = "Software"
$ws.Range("B7").Value = "Tester"
$ws.Range("C7").Value = "SDET"
$ws.Range("D7").Value = "/Users/shpendpllana/Desktop/PANDA.jpeg"
$ws.Range("E7").Value = "thebestCareer"

# ---- New row 8: NQ / Future / TRADE / lifeChanger --------------------------
$ws.Range("A8").Value = "NQ"
$ws.Range("B8").Value = "Future"
$ws.Range("C8").Value = "TRADE"
$ws.Range("D8").Value = "/Users/shpendpllana/Desktop/PANDA.jpeg"
$ws.Range("E8").Value = "lifeChanger"
$ws.Range("F8").Value = "Hum@nhrm123"
$ws.Range("G8").Value = "Hum@nhrm123"

# Pick up the hyperlink-cell look (fontId=5, xfId=0, same as F2:G7) by
# copying the existing F7:G7 formatting onto the new F8:G8 cells without
# touching the values we just wrote.
$ws.Range("F7:G7").Copy()
$ws.Range("F8:G8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Wire up the same mailto hyperlinks the other rows use.
$ws.Hyperlinks.Add($ws.Range("F8"), "mailto:Hum@nhrm123", "", "", "Hum@nhrm123")
$ws.Hyperlinks.Add($ws.Range("G8"), "mailto:Hum@nhrm123", "", "", "Hum@nhrm123")

# The active selection moved to H8 in the saved file.
$ws.Range("H8").Select()

# Cosmetic window-size tweak recorded in the workbook view.
$excel.Windows.Item(1).Width = 32000
$excel.Windows.Item(1).Height = 11940
